$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 859.6393399999999
$ws.Range("I15").Value = 859.6393399999999
$ws.Range("K15").Value = 2578.91802
$ws.Range("M15").Value = -2409.91802
$ws.Range("H33").Value = 76923624
$ws.Range("I33").Value = 76923624
$ws.Range("K33").Value = 76923624
$ws.Range("M33").Value = -76923395
$ws.Range("H40").Value = 1840
$ws.Range("I40").Value = 1740
$ws.Range("J40").Value = 1930
$ws.Range("K40").Value = 1740
$ws.Range("L40").Value = 1930
$ws.Range("M40").Value = -1565
$ws.Range("N40").Value = -2280
$ws.Range("H76").Value = 7457.095
$ws.Range("J76").Value = 8315.385
$ws.Range("L76").Value = 8315.385
$ws.Range("N76").Value = -8945.385
$ws.Range("H79").Value = 7457.095
$ws.Range("J79").Value = 8315.385
$ws.Range("L79").Value = 8315.385
$ws.Range("N79").Value = -10499.385
$ws.Range("H131").Value = 4430.625
$ws.Range("I131").Value = 542.25
$ws.Range("J131").Value = 8319
$ws.Range("K131").Value = 1626.75
$ws.Range("L131").Value = 24957
$ws.Range("M131").Value = 3413.25
$ws.Range("N131").Value = -35037
$ws.Range("H132").Value = 3312.611
$ws.Range("I132").Value = 3092.9714
$ws.Range("K132").Value = 9278.914199999999
$ws.Range("M132").Value = -6748.914199999999
$ws.Range("H138").Value = 2308.0938
$ws.Range("I138").Value = 1762.2941
$ws.Range("J138").Value = 2926.6667
$ws.Range("K138").Value = 5286.8823
$ws.Range("L138").Value = 8780.000100000001
$ws.Range("M138").Value = -146.8823000000002
$ws.Range("N138").Value = -19060.0001
$ws.Range("H141").Value = 3543.25
$ws.Range("I141").Value = 3235
$ws.Range("J141").Value = 4776.25
$ws.Range("K141").Value = 9705
$ws.Range("L141").Value = 14328.75
$ws.Range("M141").Value = -4525
$ws.Range("N141").Value = -24688.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2496.818
$ws.Range("I61").Value = 2545.1
$ws.Range("J61").Value = 2014
$ws.Range("K61").Value = 2545.1
$ws.Range("L61").Value = 2014
$ws.Range("M61").Value = -2333.1
$ws.Range("N61").Value = -2438
$ws.Range("H136").Value = 2496.818
$ws.Range("I136").Value = 2545.1
$ws.Range("J136").Value = 2014
$ws.Range("K136").Value = 7635.299999999999
$ws.Range("L136").Value = 6042
$ws.Range("M136").Value = -5085.299999999999
$ws.Range("N136").Value = -11142

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 475
$ws.Range("I22").Value = 450
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 450
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = -277
$ws.Range("N22").Value = -846

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 4501.8335
$ws.Range("I94").Value = 3062
$ws.Range("J94").Value = 5418.091
$ws.Range("K94").Value = 3062
$ws.Range("L94").Value = 5418.091
$ws.Range("M94").Value = -2611
$ws.Range("N94").Value = -6320.091

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 892.0417
$ws.Range("I114").Value = 287.9091
$ws.Range("J114").Value = 1403.2307
$ws.Range("K114").Value = 863.7273
$ws.Range("L114").Value = 4209.6921
$ws.Range("M114").Value = 2390.2727
$ws.Range("N114").Value = -10717.6921
$ws.Range("H119").Value = 7469.3335
$ws.Range("I119").Value = 4963.2
$ws.Range("J119").Value = 20000
$ws.Range("K119").Value = 14889.6
$ws.Range("L119").Value = 60000
$ws.Range("M119").Value = -10051.6
$ws.Range("N119").Value = -69676
$ws.Range("H122").Value = 1396.8695
$ws.Range("I122").Value = 1137.1333
$ws.Range("J122").Value = 1883.875
$ws.Range("K122").Value = 10234.1997
$ws.Range("L122").Value = 16954.875
$ws.Range("M122").Value = -7784.199699999999
$ws.Range("N122").Value = -21854.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 112320.11
$ws.Range("I113").Value = 200592.2
$ws.Range("J113").Value = 1980
$ws.Range("K113").Value = 200592.2
$ws.Range("L113").Value = 1980
$ws.Range("M113").Value = -198422.2
$ws.Range("N113").Value = -6320

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4159.2
$ws.Range("I16").Value = 3736.5
$ws.Range("J16").Value = 5850
$ws.Range("K16").Value = 3736.5
$ws.Range("L16").Value = 5850
$ws.Range("M16").Value = -3566.5
$ws.Range("N16").Value = -6190
$ws.Range("H40").Value = 4778.5
$ws.Range("I40").Value = 4351.7896
$ws.Range("K40").Value = 4351.7896
$ws.Range("M40").Value = -4215.7896
$ws.Range("H55").Value = 482
$ws.Range("I55").Value = 398.5
$ws.Range("J55").Value = 533.38464
$ws.Range("K55").Value = 398.5
$ws.Range("L55").Value = 533.38464
$ws.Range("M55").Value = -225.5
$ws.Range("N55").Value = -879.38464
$ws.Range("H62").Value = 19950
$ws.Range("J62").Value = 19950
$ws.Range("L62").Value = 19950
$ws.Range("N62").Value = -21198
$ws.Range("H65").Value = 19950
$ws.Range("J65").Value = 19950
$ws.Range("L65").Value = 59850
$ws.Range("N65").Value = -66090
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("H75").Value = 27086.5
$ws.Range("J75").Value = 27086.5
$ws.Range("L75").Value = 27086.5
$ws.Range("N75").Value = -28958.5
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("H78").Value = 27086.5
$ws.Range("J78").Value = 27086.5
$ws.Range("L78").Value = 81259.5
$ws.Range("N78").Value = -90619.5
$ws.Range("H80").Value = 25128
$ws.Range("J80").Value = 25128
$ws.Range("L80").Value = 25128
$ws.Range("N80").Value = -27374
$ws.Range("H82").Value = 1806.6818
$ws.Range("I82").Value = 1728.8182
$ws.Range("J82").Value = 1884.5454
$ws.Range("K82").Value = 1728.8182
$ws.Range("L82").Value = 1884.5454
$ws.Range("M82").Value = -1367.8182
$ws.Range("N82").Value = -2606.5454
$ws.Range("H83").Value = 25128
$ws.Range("J83").Value = 25128
$ws.Range("L83").Value = 75384
$ws.Range("N83").Value = -86616
$ws.Range("H85").Value = 1806.6818
$ws.Range("I85").Value = 1728.8182
$ws.Range("J85").Value = 1884.5454
$ws.Range("K85").Value = 1728.8182
$ws.Range("L85").Value = 1884.5454
$ws.Range("M85").Value = -480.8181999999999
$ws.Range("N85").Value = -4380.5454

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 1494.9
$ws.Range("I23").Value = 1105.4445
$ws.Range("J23").Value = 5000
$ws.Range("K23").Value = 1105.4445
$ws.Range("L23").Value = 5000
$ws.Range("M23").Value = -876.4445000000001
$ws.Range("N23").Value = -5458
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("H96").Value = 2966.9333
$ws.Range("I96").Value = 2750
$ws.Range("J96").Value = 3045.818
$ws.Range("K96").Value = 2750
$ws.Range("L96").Value = 3045.818
$ws.Range("M96").Value = -1377
$ws.Range("N96").Value = -5791.818
$ws.Range("H126").Value = 2668.6365
$ws.Range("I126").Value = 1900.5264
$ws.Range("J126").Value = 7533.3335
$ws.Range("K126").Value = 5701.5792
$ws.Range("L126").Value = 22600.0005
$ws.Range("M126").Value = -3231.5792
$ws.Range("N126").Value = -27540.0005
